$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.127881588408715, 0.3127903958511391, 0.8054896365839992, 8.660232485948974, 0, 9.906394106792828)
    3 = @(0.0008583669626518464, 0.0000006633126561350622, 0.8054896365839992, 0.496779210170732, 0, 1.303127877030039)
    4 = @(0.04763786555579896, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 1, 0.7443468554461139)
    5 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 0, 9.295990156953671)
    6 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $vals[$col - 2]
    }
}

$wb.Save()
